$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the post row for "「良い時間は速く過ぎる」" (row 472). All subsequent
# rows shift up by one to close the gap.
$ws.Rows.Item(472).Delete()
